# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the style of the last existing header cell (AC1, s="1")
# onto the three new header cells so they pick up the same bold/border/
# centered formatting, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-41: every player gets the same team season record.
$lastRow = 41
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 86   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 76   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
